$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dog")

# Fix typo in existing note for row 9 (trailing slash -> period)
$ws.Range("L9").Value = "Worked downhill. Did not get onto odour. Search got derailed by Koda picking up odour of a frisbee next to the search area. Probably going to exclude because ignoring a high value item is not part of the controlled evaluation protocol."

# Copy formatting from row 10 down to the two new rows (reuses existing
# date/time number-format styles instead of minting new ones)
$ws.Range("A10:L10").Copy()
$ws.Range("A11:L12").PasteSpecial(-4122)

# Add new row 11
$ws.Range("A11").Value = 45798
$ws.Range("B11").Value = "PRESENCE"
$ws.Range("C11").Value = 0.29166666666666669
$ws.Range("D11").Value = 0.51388888888888884
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 9
$ws.Range("G11").Value = "Overcast, mild"
$ws.Range("H11").Value = $true
$ws.Range("I11").Value = "2 minutes 19 seconds"
$ws.Range("J11").Value = 139
$ws.Range("K11").Value = "Primary sweeps"
$ws.Range("L11").Value = "Worked uphill. Found on the second sweep and Koda did a great job sourcing odour among complex tussock."

# Add new row 12
$ws.Range("A12").Value = 45802
$ws.Range("B12").Value = "PRESENCE"
$ws.Range("C12").Value = 0.46527777777777779
$ws.Range("D12").Value = 0.65277777777777779
$ws.Range("E12").Value = 16
$ws.Range("F12").Value = 19
$ws.Range("G12").Value = "Sunny, cool"
$ws.Range("H12").Value = $true
$ws.Range("I12").Value = "2 minutes 21 seconds"
$ws.Range("J12").Value = 141
$ws.Range("K12").Value = "Primary sweeps"
$ws.Range("L12").Value = "Worked downhill. Was on third sweep when Koda belted down hill and found it very fast, amazing - best find yet. "

# Update selection to reflect new active cell after edits
$ws.Range("A13").Select() | Out-Null
